$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "shg"
$ws.Range("C4").Value = "dfjdjg"
$ws.Range("C8").Value = "fjjgj"
$ws.Range("G8").Value = "ejeht"

$ws.Range("G8").Select()
